$wb = $excel.ActiveWorkbook

# --- Sheet "результат расчета программой" -> rename to "report" -------------
$wsReport = $wb.Worksheets.Item(2)
$wsReport.Name = "report"

# Work on the report sheet first so the other sheet ends up as the final
# active / tabSelected sheet (matches the target workbook state).
$wsReport.Activate() | Out-Null
$wsReport.Range("C23").Select() | Out-Null

# pageSetup: orientation="portrait" (dpi/copies/printer-relationship metadata
# are not reachable through the exposed PageSetup COM surface)
$wsReport.PageSetup.Orientation = 1

# --- Sheet "квартири, площі" -------------------------------------------------
$wsQuartiles = $wb.Worksheets.Item(1)
$wsQuartiles.Activate() | Out-Null

$win = $excel.ActiveWindow
$win.FreezePanes = $true

# Move the frozen pane's top-left visible cell from A92 to A2.
$win.ScrollRow = 2
$win.ScrollColumn = 1

# First select in the (unfrozen-origin) top pane, then move the real/active
# selection into the frozen bottom-left pane at E2.
$wsQuartiles.Range("K108").Select() | Out-Null
$wsQuartiles.Range("E2").Select() | Out-Null
